$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8000
$ws.Range("I62").Value = 8000
$ws.Range("K62").Value = 8000
$ws.Range("M62").Value = -7376

$ws.Range("H65").Value = 8000
$ws.Range("I65").Value = 8000
$ws.Range("K65").Value = 40000
$ws.Range("M65").Value = -36880

$ws.Range("H118").Value = 83333810
$ws.Range("J118").Value = 783.6667
$ws.Range("L118").Value = 2351.0001
$ws.Range("N118").Value = -5665.0001

$ws.Range("H132").Value = 19610808
$ws.Range("J132").Value = 2997.375
$ws.Range("L132").Value = 8992.125
$ws.Range("N132").Value = -14052.125

$ws.Range("H135").Value = 902.71875
$ws.Range("I135").Value = 564.8570999999999
$ws.Range("K135").Value = 5083.7139
$ws.Range("M135").Value = -2548.7139

$ws.Range("H137").Value = 54335.15
$ws.Range("I137").Value = 100289.72
$ws.Range("J137").Value = 2636.25
$ws.Range("K137").Value = 300869.16
$ws.Range("L137").Value = 7908.75
$ws.Range("M137").Value = -298319.16
$ws.Range("N137").Value = -13008.75

$ws.Range("H138").Value = 2895.9546
$ws.Range("J138").Value = 3721.077
$ws.Range("L138").Value = 11163.231
$ws.Range("N138").Value = -21443.231

$ws.Range("H141").Value = 11905.883
$ws.Range("I141").Value = 6300
$ws.Range("J141").Value = 53950
$ws.Range("K141").Value = 18900
$ws.Range("L141").Value = 161850
$ws.Range("M141").Value = -13720
$ws.Range("N141").Value = -172210

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7982.6753
$ws.Range("I32").Value = 4551.224
$ws.Range("J32").Value = 18457.63
$ws.Range("K32").Value = 4551.224
$ws.Range("L32").Value = 18457.63
$ws.Range("M32").Value = -4264.224
$ws.Range("N32").Value = -19031.63

$ws.Range("H45").Value = 5331857
$ws.Range("I45").Value = 10277456
$ws.Range("K45").Value = 10277456
$ws.Range("M45").Value = -10277079

$ws.Range("H63").Value = 4705.65
$ws.Range("I63").Value = 1537.5
$ws.Range("J63").Value = 7873.8
$ws.Range("K63").Value = 1537.5
$ws.Range("L63").Value = 7873.8
$ws.Range("M63").Value = -851.5
$ws.Range("N63").Value = -9245.799999999999

$ws.Range("H66").Value = 4705.65
$ws.Range("I66").Value = 1537.5
$ws.Range("J66").Value = 7873.8
$ws.Range("K66").Value = 7687.5
$ws.Range("L66").Value = 39369
$ws.Range("M66").Value = -4255.5
$ws.Range("N66").Value = -46233

$ws.Range("H132").Value = 3294
$ws.Range("I132").Value = 2065
$ws.Range("J132").Value = 6103.143
$ws.Range("K132").Value = 6195
$ws.Range("L132").Value = 18309.429
$ws.Range("M132").Value = -3665
$ws.Range("N132").Value = -23369.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 23815564
$ws.Range("I20").Value = 27783658
$ws.Range("J20").Value = 6992.5
$ws.Range("K20").Value = 27783658
$ws.Range("L20").Value = 6992.5
$ws.Range("M20").Value = -27783411
$ws.Range("N20").Value = -7486.5

$ws.Range("H138").Value = 67252.71000000001
$ws.Range("J138").Value = 67252.71000000001
$ws.Range("L138").Value = 67252.71000000001
$ws.Range("N138").Value = -77532.71000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7047.75
$ws.Range("I58").Value = 11573.8
$ws.Range("J58").Value = 3814.8572
$ws.Range("K58").Value = 11573.8
$ws.Range("L58").Value = 3814.8572
$ws.Range("M58").Value = -11370.8
$ws.Range("N58").Value = -4220.8572

$ws.Range("H122").Value = 3562.5625
$ws.Range("I122").Value = 3600.0667
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 10800.2001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -8350.2001
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 50325.31
$ws.Range("I132").Value = 31783
$ws.Range("J132").Value = 152308
$ws.Range("K132").Value = 95349
$ws.Range("L132").Value = 456924
$ws.Range("M132").Value = -92819
$ws.Range("N132").Value = -461984

$ws.Range("H136").Value = 7047.75
$ws.Range("I136").Value = 11573.8
$ws.Range("J136").Value = 3814.8572
$ws.Range("K136").Value = 34721.39999999999
$ws.Range("L136").Value = 11444.5716
$ws.Range("M136").Value = -32171.39999999999
$ws.Range("N136").Value = -16544.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 24500
$ws.Range("J63").Value = 24500
$ws.Range("L63").Value = 24500
$ws.Range("N63").Value = -25872

$ws.Range("H66").Value = 24500
$ws.Range("J66").Value = 24500
$ws.Range("L66").Value = 73500
$ws.Range("N66").Value = -80364

$ws.Range("H70").Value = 66670332
$ws.Range("I70").Value = 100002750
$ws.Range("K70").Value = 100002750
$ws.Range("M70").Value = -100002480

$ws.Range("H73").Value = 66670332
$ws.Range("I73").Value = 100002750
$ws.Range("K73").Value = 100002750
$ws.Range("M73").Value = -100001814

$ws.Range("H80").Value = 20170924
$ws.Range("I80").Value = 30848398
$ws.Range("J80").Value = 2362.3333
$ws.Range("K80").Value = 30848398
$ws.Range("L80").Value = 2362.3333
$ws.Range("M80").Value = -30847400
$ws.Range("N80").Value = -4358.3333

$ws.Range("H83").Value = 20170924
$ws.Range("I83").Value = 30848398
$ws.Range("J83").Value = 2362.3333
$ws.Range("K83").Value = 154241990
$ws.Range("L83").Value = 11811.6665
$ws.Range("M83").Value = -154236998
$ws.Range("N83").Value = -21795.6665

$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180

$ws.Range("H122").Value = 310136.2
$ws.Range("I122").Value = 406781.4
$ws.Range("K122").Value = 1220344.2
$ws.Range("M122").Value = -1217894.2

$ws.Range("H136").Value = 14424.5
$ws.Range("J136").Value = 14424.5
$ws.Range("L136").Value = 43273.5
$ws.Range("N136").Value = -48373.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 58644040
$ws.Range("I82").Value = 105558280
$ws.Range("J82").Value = 1234.5
$ws.Range("K82").Value = 105558280
$ws.Range("L82").Value = 1234.5
$ws.Range("M82").Value = -105557919
$ws.Range("N82").Value = -1956.5

$ws.Range("H85").Value = 58644040
$ws.Range("I85").Value = 105558280
$ws.Range("J85").Value = 1234.5
$ws.Range("K85").Value = 105558280
$ws.Range("L85").Value = 1234.5
$ws.Range("M85").Value = -105557032
$ws.Range("N85").Value = -3730.5

$ws.Range("H132").Value = 7447.3076
$ws.Range("I132").Value = 7683.7354
$ws.Range("J132").Value = 5839.6
$ws.Range("K132").Value = 23051.2062
$ws.Range("L132").Value = 17518.8
$ws.Range("M132").Value = -20521.2062
$ws.Range("N132").Value = -22578.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4199.2
$ws.Range("I122").Value = 2981.85
$ws.Range("J122").Value = 6633.9
$ws.Range("K122").Value = 8945.549999999999
$ws.Range("L122").Value = 19901.7
$ws.Range("M122").Value = -6495.549999999999
$ws.Range("N122").Value = -24801.7

$ws.Range("H132").Value = 28887732
$ws.Range("I132").Value = 43485150
$ws.Range("K132").Value = 130455450
$ws.Range("M132").Value = -130452920
